$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: update the source date (B5) and the CDC code (E5). The dependent
# formulas recalc automatically -- C5 (shared formula C5:C10) and F5.
$ws.Range("B5").Value = 44533
$ws.Range("E5").Value = 7534
$ws.Range("F5").Formula = "=IF(E5 > 0, DATE(2001,5,1)+E5-1, ""cdc inválido"")"

# Move the active cell/selection from H10 to J9.
$ws.Range("J9").Select()
